$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.812.85'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.650.03'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').Value = "'217.01"
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('D6').Value = "'0.505"
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = "'19.33"
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('D11').Value = "'0.0845"
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '1.873.97'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').Value = '1.677.53'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = "'4.23"
$ws.Range('E14').Value = '  +1.85%  '
$ws.Range('D15').Value = "'0.534"
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').Value = "'65.86"
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '26.805.35'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('E18').Value = '  -0.09%  '
$ws.Range('D19').Value = "'217.23"
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('E20').Value = '  +0.70%  '
$ws.Range('E21').Value = '  +0.49%  '
$ws.Range('D22').Value = "'2.43"
$ws.Range('E22').Value = '  +15.33%  '
$ws.Range('D23').Value = "'6.34"
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('D25').Value = "'145.90"
$ws.Range('E25').Value = '  -1.16%  '
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('E28').Value = '  +3.88%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').Value = "'0.0522"
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('E32').Value = '  -0.72%  '
$ws.Range('E33').Value = '  +0.82%  '
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').Value = '1.279.09'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  +1.88%  '
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').Value = "'0.544"
$ws.Range('E38').Value = '  +5.22%  '
$ws.Range('D39').Value = "'0.836"
$ws.Range('E39').Value = '  +3.09%  '
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').Value = "'0.821"
$ws.Range('E41').Value = '  +1.57%  '
$ws.Range('D42').Value = "'2.24"
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('D44').Value = '1.799.24'
$ws.Range('D45').Value = "'92.27"
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('E46').Value = '  +6.36%  '
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('E49').Value = '  +0.58%  '
$ws.Range('D50').Value = "'7.81"
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').Value = "'0.0986"
$ws.Range('E51').Value = '  +1.87%  '
